$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56 (pushes existing rows 56..80 down to 57..81)
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record
$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 44875
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = 100112013
$ws.Cells.Item(56, 7).Value = "Alcachofa"
$ws.Cells.Item(56, 8).Value = "Española"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 60
$ws.Cells.Item(56, 11).Value = 10000
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 10000
$ws.Cells.Item(56, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(56, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(56, 16).Value = 333
$ws.Cells.Item(56, 17).Value = 30
$ws.Cells.Item(56, 18).Value = "Hortaliza"
